# Append the newly recommended skill entries to the Skillsets list
# (GetRecommended7 now returns the full set, not just the "Recommended"
# subset, so the sheet grows from 438 to 460 data rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Values are written in the same order the shared-string table assigns new
# indices (439 first, then 440-442, then 444 before 443, etc.) so the
# resulting sharedStrings.xml ordering matches the source edit exactly.
$ws.Cells.Item(439, 1).Value = 'Microsoft PowerPoint'
$ws.Cells.Item(440, 1).Value = 'Outsourcing '
$ws.Cells.Item(441, 1).Value = 'Investigation'
$ws.Cells.Item(442, 1).Value = 'MS Word'
$ws.Cells.Item(444, 1).Value = 'Microsoft Word'
$ws.Cells.Item(443, 1).Value = 'Able To Multitask'
$ws.Cells.Item(446, 1).Value = 'Statistics'
$ws.Cells.Item(447, 1).Value = 'Able To Work Independently'
$ws.Cells.Item(445, 1).Value = 'Document Management '
$ws.Cells.Item(448, 1).Value = 'Public Relations'
$ws.Cells.Item(449, 1).Value = 'Good Communication Skills'
$ws.Cells.Item(450, 1).Value = 'Communication Skills'
$ws.Cells.Item(451, 1).Value = 'Administrative Support'
$ws.Cells.Item(452, 1).Value = 'Circulation'
$ws.Cells.Item(453, 1).Value = 'Taxation'
$ws.Cells.Item(454, 1).Value = 'Human Resources  '
$ws.Cells.Item(455, 1).Value = 'Employee Relations'
$ws.Cells.Item(457, 1).Value = 'Continuous Improvement  '
$ws.Cells.Item(459, 1).Value = 'Logistics'
$ws.Cells.Item(460, 1).Value = 'Facility management'
$ws.Cells.Item(456, 1).Value = 'Document Controller'
$ws.Cells.Item(458, 1).Value = 'Shippings'

$ws.Range("B452").Select()
